$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.689.05"
$ws.Range("E2").Value = "  -3.16%  "

$ws.Range("D3").Value = "2.096.64"
$ws.Range("E3").Value = "  -1.37%  "

$ws.Range("E4").Value = "  -0.32%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "343.54"
$c.Style = "Normal"

$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("E7").Value = "  -2.46%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.4404"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -3.17%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "53.04"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.80%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.09202"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.18%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.171"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.04%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "24.92"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("D13").Value = "2.099.06"
$ws.Range("E13").Value = "  -1.90%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.760"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.53%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "8.193"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.84%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "99.45"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -2.96%  "

$ws.Range("E17").Value = "  -2.36%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.45%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "20.74"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +6.46%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.06636"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.15%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.38%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.190"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.61%  "

$ws.Range("D23").Value = "29.756.49"
$ws.Range("E23").Value = "  -3.20%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "12.59"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.54%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.321"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.03%  "

$ws.Range("D26").Value = "2.343.12"
$ws.Range("E26").Value = "  -1.83%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "21.88"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.96%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "162.49"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.39%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.527"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.57%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "132.44"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.90%  "

$ws.Range("E31").Value = "  -5.92%  "

$ws.Range("E32").Value = "  -3.34%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.648"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.38%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.157"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.64%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.946"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.84%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "6.040"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.18%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "10.43"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.03%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02563"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.52%  "

$ws.Range("E39").Value = "  -2.68%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "12.41"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.32%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.2236"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.84%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.6866"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.83%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.293"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.13%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.6651"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.90%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "14.13"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -4.38%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.293"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.15%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.618"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.91%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.00000000349"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -4.85%  "

$ws.Range("E49").Value = "  -3.20%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "81.89"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.33%  "

$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.3265"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.90%  "
